$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Z (26th column) to hold the new "STAT" field.
# This shifts all existing columns Z..AJ one place to the right (-> AA..AK)
# exactly as Excel's native column insert does, keeping header/value/style mapping intact.
$ws.Columns("Z:Z").Insert()

# Populate the newly inserted STAT column
$ws.Range("Z1").Value = "STAT"
$ws.Range("Z2").Value = "A"
$ws.Range("Z3").Value = "S"

# Match the target column width (stored width="8" in the XML)
$ws.Columns("Z:Z").ColumnWidth = 7.15

# Reflect final selection/scroll state recorded in the saved workbook
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 26
$ws.Range("AJ10").Select()
